# Update Name of Algo
# Apply updated RandomForest imputation results to the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 6.048599999999989
$ws.Range("A9").Value  = -20.37649999999997
$ws.Range("B11").Value = 5.490799999999996
$ws.Range("A18").Value = -22.97790000000001
$ws.Range("A20").Value = -22.11690000000002
$ws.Range("C21").Value = -13.34820000000001
